$d = $word.ActiveDocument

$br = [char]11

function Replace-ExactText($findText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $findText"
    }
    $rng.Text = $newText
}

# --- 1) "Programa" section, Portuguese paragraph ---
$ptOld = "1) Carga e Força elétrica: carga elétrica; condutores e isolantes; lei de Coulomb; quantização e conservação de cargas.2) Campo Elétrico: conceito; linhas de campo; carga pontual e dipolo elétrico, distribuição contínua.3) A Lei de Gauss: fluxo; aplicações em simetrias cilíndricas, planares e esféricas.4) Potencial Elétrico: conceito e cálculo; energia, potencial e campo elétrico, superfícies equipotenciais; carga puntiforme, dipolo elétrico e distribuições contínuas.5) Capacitores e Dielétricos: capacitância, energia e cálculo; associações, dielétrico.6) Corrente e Resistência Elétrica: corrente e densidade, resistência, Resistividade e Condutividade em função da temperatura; lei de Ohm, potência, semicondutores e supercondutores.7) Campos Magnéticos: lei de Biot-Savart.8) Lei de Ampère e aplicações; campo magnético de uma espira, solenoide e toroides.9) Indução Eletromagnética: conceitos; Lei de indução de Faraday; Lei de Lenz;10) Propriedades magnéticas da matéria;11) Equações de Maxwell."

$ptNew = "1) Carga e Força elétrica: carga elétrica; condutores e isolantes; lei de Coulomb; quantização e conservação de cargas." + $br + `
"2) Campo Elétrico: conceito; linhas de campo; carga pontual e dipolo elétrico, distribuição contínua." + $br + `
"3) A Lei de Gauss: fluxo; aplicações em simetrias cilíndricas, planares e esféricas." + $br + `
"4) Potencial Elétrico: conceito e cálculo; energia, potencial e campo elétrico, superfícies equipotenciais; carga puntiforme, dipolo elétrico e distribuições contínuas." + $br + `
"5) Capacitores e Dielétricos: capacitância, energia e cálculo; associações, dielétrico." + $br + `
"6) Corrente e Resistência Elétrica: corrente e densidade, resistência, Resistividade e Condutividade em função da temperatura; lei de Ohm, potência, semicondutores e supercondutores." + $br + `
"7) Campos Magnéticos: lei de Biot-Savart." + $br + `
"8) Lei de Ampère e aplicações; campo magnético de uma espira, solenoide e toroides." + $br + `
"9) Indução Eletromagnética: conceitos; Lei de indução de Faraday; Lei de Lenz;" + $br + `
"10) Propriedades magnéticas da matéria;" + $br + `
"11) Equações de Maxwell."

Replace-ExactText $ptOld $ptNew

# --- 2) "Programa" section, English (italic) paragraph ---
$enOld = "1) Electric charge and electric force: electric charge; conductors and insulators; Coulomb's law; quantization and conservation.2) Electric field: concepts; field lines; point charge and dipole, continuous distribution.3) Gauss' law: flow; applications in cylindrical, flat and spherical geometries.4) Electric potential: concept and calculation; energy, potential and electric field, equipotential surfaces; punctual loads, electric dipole and continuous distributions.5) Capacitors and dielectrics: capacitance, energy and calculation, associations, dielectrics.6) Electric current and resistance: current density, resistance and resistivity as a function of temperature; Ohm's law, power, semiconductors and superconductors.7) Magnetic field: Biot-Savart law.8) Ampère's law and applications: magnetic field of a coil, solenoid, and toroids.9) Electromagnetic induction and inductance: Faraday's law, Lenz's law.10) Magnetic properties of matter.11) Maxwell's equations."

$enNew = "1) Electric charge and electric force: electric charge; conductors and insulators; Coulomb's law; quantization and conservation." + $br + `
"2) Electric field: concepts; field lines; point charge and dipole, continuous distribution." + $br + `
"3) Gauss' law: flow; applications in cylindrical, flat and spherical geometries." + $br + `
"4) Electric potential: concept and calculation; energy, potential and electric field, equipotential surfaces; punctual loads, electric dipole and continuous distributions." + $br + `
"5) Capacitors and dielectrics: capacitance, energy and calculation, associations, dielectrics." + $br + `
"6) Electric current and resistance: current density, resistance and resistivity as a function of temperature; Ohm's law, power, semiconductors and superconductors." + $br + `
"7) Magnetic field: Biot-Savart law." + $br + `
"8) Ampère's law and applications: magnetic field of a coil, solenoid, and toroids." + $br + `
"9) Electromagnetic induction and inductance: Faraday's law, Lenz's law." + $br + `
"10) Magnetic properties of matter." + $br + `
"11) Maxwell's equations."

Replace-ExactText $enOld $enNew

# --- 3) "Bibliografia" section ---
$bibOld = "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.3, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.3, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 3, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008)."

$bibNew = "NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 3, Edgard Blucher (2008)." + $br + `
"RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.3, LTC (2008)." + $br + `
"TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.3, LTC (2008)." + $br + `
"SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 3, Pearson Addison Wesley (2009)." + $br + `
"JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 3, Thomson Pioneira (2008)."

Replace-ExactText $bibOld $bibNew
